# Update "want to go" (想去人数) counts in F column on both the "展览"
# and "全部类型" sheets to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 281
    3  = 1442
    4  = 174
    10 = 140
    11 = 4
    12 = 4777
    14 = 7052
    15 = 43
    16 = 61
    20 = 7
    21 = 4188
    22 = 1285
    23 = 84
    25 = 2768
    27 = 556
    29 = 393
    30 = 394
    31 = 415
    33 = 55
    34 = 1655
    35 = 1065
    36 = 74
    37 = 890
    38 = 90
    43 = 24
    45 = 736
}

$updates4 = @{
    2  = 281
    3  = 1442
    4  = 174
    10 = 140
    11 = 4
    12 = 4777
    14 = 7052
    15 = 43
    16 = 61
    20 = 7
    21 = 4188
    22 = 1286
    23 = 84
    25 = 2768
    27 = 556
    29 = 393
    30 = 394
    31 = 415
    33 = 55
    34 = 1655
    35 = 1065
    36 = 74
    37 = 890
    38 = 90
    43 = 24
    45 = 736
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
